# Update the MODS wrapper elements used to build the Fedora "object" XML
# snippet for each row: the old workflow wrapped the MODS payload in
#   <update type="MODS"> ... </update>
# the new workflow wraps it in
#   <datastream type="md_descriptive" operation="update"> ... </datastream>
#
# These two fragments live in the shared-string table and are referenced
# by the header cells C1 (wrapper open tag, following the pid attribute)
# and BM1 (wrapper close tag, preceding the closing </object> tag).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value2 = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("BM1").Value2 = "</mods:mods></datastream></object>"

# The author's last selection before saving moved from J1 to BM1.
$ws.Range("BM1").Select()
